# draft-gandhi-mpls-ioam-sr-02.pptx -- "Add files via upload"
#
# The deck's "History of the Draft" slide drops the two bullet
# paragraphs about the IETF 108 MPLS WG presentation ("Jul 2020" /
# "Presented in draft-gandhi-mpls-ioam-sr-02 at IETF 108 in MPLS WG"),
# while leaving the preceding "...at IETF 107 in MPLS WG Interim"
# bullet (and the trailing blank paragraph) untouched.

$p = $ppt.ActivePresentation

# Locate the slide/shape that carries the draft history bullets
# instead of hard-coding indices, so the script is resilient to any
# reordering.
$targetSlide = $null
$targetShape = $null

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
        $shape = $slide.Shapes.Item($shi)
        if ($shape.HasTextFrame -and $shape.TextFrame.HasText) {
            if ($shape.TextFrame.TextRange.Text -like "*at IETF 108 in MPLS WG*") {
                $targetSlide = $slide
                $targetShape = $shape
            }
        }
    }
}

if ($targetShape -ne $null) {
    $tr = $targetShape.TextFrame.TextRange
    $count = $tr.Paragraphs().Count

    # Find the two consecutive paragraphs to remove:
    #   "Jul 2020"
    #   "Presented in draft-gandhi-mpls-ioam-sr-02 at IETF 108 in MPLS WG"
    $firstIdx = -1
    for ($i = 1; $i -le $count; $i++) {
        $paraText = $tr.Paragraphs($i, 1).Text
        if ($paraText -like "Jul 2020*") {
            if (($i + 1) -le $count) {
                $nextText = $tr.Paragraphs($i + 1, 1).Text
                if ($nextText -like "*Presented in*" -and $nextText -like "*at IETF 108 in MPLS WG*") {
                    $firstIdx = $i
                }
            }
        }
    }

    if ($firstIdx -ne -1) {
        $para1 = $tr.Paragraphs($firstIdx, 1)
        $para2 = $tr.Paragraphs($firstIdx + 1, 1)
        $combinedLength = $para1.Length + $para2.Length
        $toRemove = $tr.Characters($para1.Start, $combinedLength)
        $toRemove.Delete()
    }
}
